# Applies the odds updates described in the commit "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("L3").Value = 6.5
$ws.Range("AH3").Value = 11

# Row 4
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5

# Row 9
$ws.Range("G9").Value = 1.44
$ws.Range("I9").Value = 6.25
$ws.Range("Q9").Value = 1.5
$ws.Range("R9").Value = 2.5
$ws.Range("AD9").Value = 9
$ws.Range("AH9").Value = 21
$ws.Range("AJ9").Value = 19
$ws.Range("BC9").Value = 126
$ws.Range("BD9").Value = 151

# Row 10
$ws.Range("BD10").Value = 126
